# Update cryptocurrency price/volume figures as scraped on Thu Sep  7 18:42:58 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '25.766.82' }
    @{ Cell = 'E2'; Value = '  +0.09%  ' }
    @{ Cell = 'D3'; Value = '1.628.42' }
    @{ Cell = 'E3'; Value = '  -0.43%  ' }
    @{ Cell = 'D4'; Value = '0.996' }
    @{ Cell = 'E4'; Value = '  -0.70%  ' }
    @{ Cell = 'D5'; Value = '213.77' }
    @{ Cell = 'E6'; Value = '  +0.04%  ' }
    @{ Cell = 'E7'; Value = '  -0.77%  ' }
    @{ Cell = 'E8'; Value = '  -1.16%  ' }
    @{ Cell = 'E9'; Value = '  -1.02%  ' }
    @{ Cell = 'E10'; Value = '  +0.09%  ' }
    @{ Cell = 'D11'; Value = '0.0788' }
    @{ Cell = 'E12'; Value = '  -0.07%  ' }
    @{ Cell = 'D13'; Value = '1.853.66' }
    @{ Cell = 'E13'; Value = '  -0.33%  ' }
    @{ Cell = 'D14'; Value = '1.624.45' }
    @{ Cell = 'E14'; Value = '  -0.85%  ' }
    @{ Cell = 'D15'; Value = '0.550' }
    @{ Cell = 'E15'; Value = '  -0.76%  ' }
    @{ Cell = 'E16'; Value = '  -1.17%  ' }
    @{ Cell = 'D17'; Value = '62.60' }
    @{ Cell = 'E17'; Value = '  -0.21%  ' }
    @{ Cell = 'D18'; Value = '25.779.40' }
    @{ Cell = 'E18'; Value = '  +0.07%  ' }
    @{ Cell = 'E19'; Value = '  -0.74%  ' }
    @{ Cell = 'E20'; Value = '  -0.36%  ' }
    @{ Cell = 'D21'; Value = '190.57' }
    @{ Cell = 'E21'; Value = '  -1.42%  ' }
    @{ Cell = 'D22'; Value = '9.89' }
    @{ Cell = 'E22'; Value = '  -0.59%  ' }
    @{ Cell = 'D23'; Value = '6.27' }
    @{ Cell = 'E23'; Value = '  -0.03%  ' }
    @{ Cell = 'E24'; Value = '  -0.77%  ' }
    @{ Cell = 'E25'; Value = '  +0.25%  ' }
    @{ Cell = 'D26'; Value = '142.15' }
    @{ Cell = 'E26'; Value = '  +1.50%  ' }
    @{ Cell = 'E27'; Value = '  +1.04%  ' }
    @{ Cell = 'E28'; Value = '  -0.56%  ' }
    @{ Cell = 'D29'; Value = '15.47' }
    @{ Cell = 'E29'; Value = '  -0.15%  ' }
    @{ Cell = 'E30'; Value = '  -0.93%  ' }
    @{ Cell = 'E31'; Value = '  -0.14%  ' }
    @{ Cell = 'E32'; Value = '  -0.53%  ' }
    @{ Cell = 'E33'; Value = '  -0.97%  ' }
    @{ Cell = 'E34'; Value = '  -0.48%  ' }
    @{ Cell = 'E35'; Value = '  -0.17%  ' }
    @{ Cell = 'D36'; Value = '0.902' }
    @{ Cell = 'E36'; Value = '  +0.27%  ' }
    @{ Cell = 'D37'; Value = '1.140.56' }
    @{ Cell = 'E37'; Value = '  +2.05%  ' }
    @{ Cell = 'E38'; Value = '  -0.84%  ' }
    @{ Cell = 'E40'; Value = '  -0.48%  ' }
    @{ Cell = 'E42'; Value = '  -1.67%  ' }
    @{ Cell = 'E43'; Value = '  +0.62%  ' }
    @{ Cell = 'D44'; Value = '100.31' }
    @{ Cell = 'E44'; Value = '  +0.58%  ' }
    @{ Cell = 'D45'; Value = '0.798' }
    @{ Cell = 'E45'; Value = '  -0.27%  ' }
    @{ Cell = 'D46'; Value = '1.765.75' }
    @{ Cell = 'E46'; Value = '  -0.20%  ' }
    @{ Cell = 'E47'; Value = '  +0.61%  ' }
    @{ Cell = 'D48'; Value = '55.28' }
    @{ Cell = 'E48'; Value = '  +0.15%  ' }
    @{ Cell = 'E49'; Value = '  +6.14%  ' }
    @{ Cell = 'E50'; Value = '  +2.05%  ' }
    @{ Cell = 'D51'; Value = '0.416' }
    @{ Cell = 'E51'; Value = '  -0.43%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (e.g. "213.77") keep their
    # exact original formatting instead of being parsed into floating point numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
